$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '24.836.29'
$ws.Cells.Item(2, 5).Value = '  +0.76%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.705.17'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.19%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '316.70'
$ws.Cells.Item(5, 5).Value = '  -0.53%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  +0.20%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.3940'
$ws.Cells.Item(7, 5).Value = '  -0.29%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.4053'
$ws.Cells.Item(8, 5).Value = '  +0.25%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '1.524'
$ws.Cells.Item(9, 5).Value = '  -1.26%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '1.002'
$ws.Cells.Item(10, 5).Value = '  +0.20%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '53.38'
$ws.Cells.Item(11, 5).Value = '  -2.21%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.08912'
$ws.Cells.Item(12, 5).Value = '  +0.93%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '7.328'
$ws.Cells.Item(13, 5).Value = '  +0.42%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '23.91'
$ws.Cells.Item(14, 5).Value = '  +1.86%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '8.016'
$ws.Cells.Item(15, 5).Value = '  +4.52%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.00001329'
$ws.Cells.Item(16, 5).Value = '  -0.28%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '1.708.59'
$ws.Cells.Item(17, 5).Value = '  +0.73%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '100.54'
$ws.Cells.Item(18, 5).Value = '  -0.87%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.07047'
$ws.Cells.Item(19, 5).Value = '  -0.74%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '19.78'
$ws.Cells.Item(20, 5).Value = '  -0.23%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '7.107'
$ws.Cells.Item(21, 5).Value = '  +2.86%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.21%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '14.53'
$ws.Cells.Item(23, 5).Value = '  +2.46%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '24.823.77'
$ws.Cells.Item(24, 5).Value = '  +0.78%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.241'
$ws.Cells.Item(25, 5).Value = '  +2.19%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.361'
$ws.Cells.Item(26, 5).Value = '  +1.15%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '22.84'
$ws.Cells.Item(27, 5).Value = '  +1.61%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '162.49'
$ws.Cells.Item(28, 5).Value = '  +1.64%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '8.474'
$ws.Cells.Item(29, 5).Value = '  +11.19%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '136.90'
$ws.Cells.Item(30, 5).Value = '  +1.83%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '5.180'

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '7.679'
$ws.Cells.Item(32, 5).Value = '  +3.16%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.08899'
$ws.Cells.Item(33, 5).Value = '  +3.51%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.087'

# Row 35
$ws.Cells.Item(35, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.987'
$ws.Cells.Item(35, 5).Value = '  +2.17%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'FraxShare'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '11.18'
$ws.Cells.Item(36, 5).Value = '  -2.93%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.2766'
$ws.Cells.Item(37, 5).Value = '  +0.25%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '14.62'
$ws.Cells.Item(38, 5).Value = '  -1.33%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.02789'
$ws.Cells.Item(39, 5).Value = '  -0.51%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Stellar'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.09213'
$ws.Cells.Item(40, 5).Value = '  +1.01%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.461'
$ws.Cells.Item(41, 5).Value = '  -0.30%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.7732'
$ws.Cells.Item(42, 5).Value = '  -0.47%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Decentraland'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.7222'
$ws.Cells.Item(43, 5).Value = '  -1.00%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '15.79'
$ws.Cells.Item(44, 5).Value = '  +1.32%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.582'
$ws.Cells.Item(45, 5).Value = '  +2.51%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.35%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.003'
$ws.Cells.Item(47, 5).Value = '  +0.22%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '140.74'
$ws.Cells.Item(48, 5).Value = '  -1.05%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.325'
$ws.Cells.Item(49, 5).Value = '  -4.36%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '91.10'
$ws.Cells.Item(50, 5).Value = '  +2.84%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.08006'
$ws.Cells.Item(51, 5).Value = '  -0.55%  '
